# Update "展览" (sheet 1) and "全部类型" (sheet 4) worksheets:
#  - bump several "want to go" (F column) counts
#  - append a new row (31) for the 上饶·ETI动漫节 event
# "演出" (sheet 2) and "本地生活" (sheet 3) are left untouched.

$wb = $excel.ActiveWorkbook

foreach ($sheetIndex in 1,4) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # --- Updated F-column ("想去人数") counts for existing rows ---
    $ws.Range("F2").Value = 1960
    $ws.Range("F5").Value = 40
    $ws.Range("F7").Value = 1639
    $ws.Range("F9").Value = 646
    $ws.Range("F10").Value = 366
    $ws.Range("F11").Value = 102
    $ws.Range("F16").Value = 139
    $ws.Range("F17").Value = 114
    $ws.Range("F18").Value = 130
    $ws.Range("F19").Value = 3801
    $ws.Range("F20").Value = 8
    $ws.Range("F21").Value = 19
    $ws.Range("F22").Value = 432
    $ws.Range("F24").Value = 703
    $ws.Range("F25").Value = 446
    $ws.Range("F28").Value = 1588
    $ws.Range("F29").Value = 18
    $ws.Range("F30").Value = 154

    # Sheet 4 ("全部类型") has one extra change not present on sheet 1
    if ($sheetIndex -eq 4) {
        $ws.Range("F27").Value = 0
    }

    # --- Append new row 31 ---
    # Copy formatting from A30 (bold/centered/bordered "index" style) to A31,
    # then overwrite with the correct index value.
    $ws.Range("A30").Copy($ws.Range("A31"))
    $ws.Range("A31").Value = 30

    # B31 holds a literal text date string (like the other rows), not a real
    # date value, so force text with a leading apostrophe and then clear the
    # formatting Excel auto-applies so no stray number format sticks around.
    $ws.Range("B31").Value = "'2024-06-10"
    $ws.Range("B31").ClearFormats()

    $ws.Range("C31").Value = "上饶·ETI动漫节"
    $ws.Range("D31").Value = "滨江东路与体育馆路交叉口西100米 力加体育综合运动中心"
    $ws.Range("E31").Value = "2024.06.10 10:00-06.10 16:00"

    # F31 differs between the two sheets
    if ($sheetIndex -eq 4) {
        $ws.Range("F31").Value = 2
    } else {
        $ws.Range("F31").Value = 1
    }

    $ws.Range("G31").Value = 36.6
    $ws.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=83422"
    $ws.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202403/vvJKFJal1711460768984.jpeg"
}
